# FHIR-55597: remove the "US Core PMO ServiceRequest Profile" row that was
# left over from a partially-applied change. Deleting the entire row shifts
# every row below it up by one, which matches the diff (row 43 disappears,
# former rows 44-57 become rows 43-56, content/values unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("43").Delete()
